$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(2, 6, 3, 3, 5, 6, 3, 5, 3, 3, 7, 6, 5, 3, 4, 3, 6, 6, 3, 1, 2)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
